$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text: "Avaliacao" -> "Avaliação"
$ws.Range("F1").Value = "Avaliação"

# Fill in previously-blank "Grupo"/UF numeric values
$ws.Range("E6").Value = 2
$ws.Range("E11").Value = 5
$ws.Range("E12").Value = 99
$ws.Range("E15").Value = 5
$ws.Range("E21").Value = 5
$ws.Range("E22").Value = 99
$ws.Range("E32").Value = 99
$ws.Range("E35").Value = 5
$ws.Range("E38").Value = 99
$ws.Range("E45").Value = 99
$ws.Range("E46").Value = 99
$ws.Range("E48").Value = 5
$ws.Range("E61").Value = 99
$ws.Range("E65").Value = 99
$ws.Range("E66").Value = 99
$ws.Range("E70").Value = 5
$ws.Range("E71").Value = 5
$ws.Range("E72").Value = 5
$ws.Range("E73").Value = 5
$ws.Range("E74").Value = 99
$ws.Range("G82").Value = 99
$ws.Range("E88").Value = 99
$ws.Range("E89").Value = 2
$ws.Range("E90").Value = 2
$ws.Range("E91").Value = 99
$ws.Range("E113").Value = 99
$ws.Range("E115").Value = 99
$ws.Range("E116").Value = 5
$ws.Range("E130").Value = 99
$ws.Range("E136").Value = 99
$ws.Range("E137").Value = 99
$ws.Range("E154").Value = 6
$ws.Range("E168").Value = 2
$ws.Range("E177").Value = 99
$ws.Range("E180").Value = 99
$ws.Range("E183").Value = 99
$ws.Range("E184").Value = 6
$ws.Range("E204").Value = 6
$ws.Range("E209").Value = 2
$ws.Range("G218").Value = 99
$ws.Range("G222").Value = 99
$ws.Range("E226").Value = 2
$ws.Range("E227").Value = 2
$ws.Range("E228").Value = 2
$ws.Range("E231").Value = 99
